$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 35, shifting the existing rows 35-45 down to 36-46.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new weekly record.
$ws.Range("A35").Value = 1
$ws.Range("B35").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C35").Value = "Arica y Parinacota"
$ws.Range("D35").Value = 44825
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 100112052
$ws.Range("G35").Value = "Albahaca"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 1800
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 1900
$ws.Range("N35").Value = "$/paquete"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 1900
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"
